$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test 2")

# New "LUNES" tags in column I for rows 2-5
$ws.Range("I2").Value = "LUNES"
$ws.Range("I3").Value = "LUNES"
$ws.Range("I4").Value = "LUNES"
$ws.Range("I5").Value = "LUNES"

# Row 4 updates
$ws.Range("A4").Value = "PRIMERO"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.4

# Row 5 updates
$ws.Range("B5").Value = "p3"
$ws.Range("C5").Value = "o3"
$ws.Range("D5").Value = 5.4
$ws.Range("E5").Value = 0.4

# Row 6: clear out A6:E6 (drops the "TERCERO" / p4 / o4 entry entirely)
$ws.Range("A6:E6").ClearContents()

# Update selection to match the recorded cursor position
[void]$ws.Range("G23").Select()
